$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.030.89"
$ws.Range("E2").Value = "  -0.99%  "

$ws.Range("D3").Value = "1.826.28"
$ws.Range("E3").Value = "  -0.34%  "

$ws.Range("D4").Value = "'1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.45%  "

$ws.Range("D5").Value = "'310.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.48%  "

$ws.Range("D6").Value = "'1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.38%  "

$ws.Range("D7").Value = "'0.4639"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.08%  "

$ws.Range("D8").Value = "'0.3655"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.89%  "

$ws.Range("D9").Value = "'0.07239"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.74%  "

$ws.Range("D10").Value = "'0.8601"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.83%  "

$ws.Range("D11").Value = "'19.88"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.84%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.910.28"
$ws.Range("E12").Value = "  +1.67%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.07794"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.36%  "

$ws.Range("D14").Value = "'5.330"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.90%  "

$ws.Range("D15").Value = "'91.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.18%  "

$ws.Range("D16").Value = "'6.497"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.08%  "

$ws.Range("D17").Value = "'1.007"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.26%  "

$ws.Range("D18").Value = "'0.000008661"
$ws.Range("D18").Style = "Normal"

$ws.Range("D20").Value = "'14.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.15%  "

$ws.Range("D21").Value = "26.814.49"
$ws.Range("E21").Value = "  -2.59%  "

$ws.Range("D22").Value = "'5.152"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.52%  "

$ws.Range("E23").Value = "  -1.25%  "

$ws.Range("D24").Value = "2.132.90"
$ws.Range("E24").Value = "  +1.88%  "

$ws.Range("D25").Value = "'151.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("D26").Value = "'1.840"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.83%  "

$ws.Range("D27").Value = "'18.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.60%  "

$ws.Range("D28").Value = "'2.065"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.89%  "

$ws.Range("D29").Value = "'5.109"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.26%  "

$ws.Range("D30").Value = "'115.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.52%  "

$ws.Range("E31").Value = "  -1.75%  "

$ws.Range("D32").Value = "'2.952"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.21%  "

$ws.Range("D33").Value = "'4.421"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.71%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.131"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.69%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7204"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.88%  "

$ws.Range("E36").Value = "  -2.04%  "

$ws.Range("D37").Value = "'0.05237"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.93%  "

$ws.Range("D38").Value = "'2.427"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.92%  "

$ws.Range("E39").Value = "  -1.57%  "

$ws.Range("D40").Value = "'2.937"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.16%  "

$ws.Range("D41").Value = "'7.159"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.91%  "

$ws.Range("D42").Value = "'0.5161"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.76%  "

$ws.Range("D43").Value = "'0.1627"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.99%  "

$ws.Range("D44").Value = "'0.8588"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -15.08%  "

$ws.Range("D45").Value = "'8.184"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.42%  "

$ws.Range("D46").Value = "'0.4795"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.90%  "

$ws.Range("E47").Value = "  -0.46%  "

$ws.Range("D48").Value = "'10.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.12%  "

$ws.Range("D49").Value = "'102.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.26%  "

$ws.Range("D50").Value = "'0.06236"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.03%  "
